$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression
$ws.Range("B2").Value = 3973184152826392
$ws.Range("C2").Value = 3973184152826392
$ws.Range("D2").Value = 3973184152826392

# Row 3 - RandomForestRegressor
$ws.Range("B3").Value = 2718773608407.162
$ws.Range("C3").Value = 2841518016381.994
$ws.Range("D3").Value = 119384728651386.7

# Row 4 - label change: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 2716240479202.796
$ws.Range("C4").Value = 2569039485745.833
$ws.Range("D4").Value = 37125781725473.19

# Row 5 - label change: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 150891519516656.6
$ws.Range("C5").Value = 211589390198999.9
$ws.Range("D5").Value = 2111364979239140
